# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the Alcachofa data block
# (rows 609-610), shifting the existing rows 609-664 down to 611-666.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("609:610").Insert()

# Row 609: Española / Extra
$ws.Cells.Item(609, 1).Value2 = 3
$ws.Cells.Item(609, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(609, 3).Value2 = "Coquimbo"
$ws.Cells.Item(609, 4).Value2 = 45223
$ws.Cells.Item(609, 5).Value2 = 5
$ws.Cells.Item(609, 6).Value2 = 100112013
$ws.Cells.Item(609, 7).Value2 = "Alcachofa"
$ws.Cells.Item(609, 8).Value2 = "Española"
$ws.Cells.Item(609, 9).Value2 = "Extra"
$ws.Cells.Item(609, 10).Value2 = 12000
$ws.Cells.Item(609, 11).Value2 = 350
$ws.Cells.Item(609, 12).Value2 = 350
$ws.Cells.Item(609, 13).Value2 = 350
$ws.Cells.Item(609, 14).Value2 = "$/unidad"
$ws.Cells.Item(609, 15).Value2 = "Llay Llay"
$ws.Cells.Item(609, 16).Value2 = 350
$ws.Cells.Item(609, 17).Value2 = 1
$ws.Cells.Item(609, 18).Value2 = "Hortaliza"

# Row 610: Española / Primera
$ws.Cells.Item(610, 1).Value2 = 3
$ws.Cells.Item(610, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(610, 3).Value2 = "Coquimbo"
$ws.Cells.Item(610, 4).Value2 = 45223
$ws.Cells.Item(610, 5).Value2 = 5
$ws.Cells.Item(610, 6).Value2 = 100112013
$ws.Cells.Item(610, 7).Value2 = "Alcachofa"
$ws.Cells.Item(610, 8).Value2 = "Española"
$ws.Cells.Item(610, 9).Value2 = "Primera"
$ws.Cells.Item(610, 10).Value2 = 13000
$ws.Cells.Item(610, 11).Value2 = 280
$ws.Cells.Item(610, 12).Value2 = 280
$ws.Cells.Item(610, 13).Value2 = 280
$ws.Cells.Item(610, 14).Value2 = "$/unidad"
$ws.Cells.Item(610, 15).Value2 = "Llay Llay"
$ws.Cells.Item(610, 16).Value2 = 280
$ws.Cells.Item(610, 17).Value2 = 1
$ws.Cells.Item(610, 18).Value2 = "Hortaliza"
